$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 473, pushing the existing rows 473:498
# (and everything below) down to 475:500.
$ws.Rows("473:474").Insert()

# Populate the newly inserted row 473 (Primera) with this week's data.
$ws.Range("A473").Value = 11
$ws.Range("B473").Value = "Vega Monumental Concepción"
$ws.Range("C473").Value = "Bíobío"
$ws.Range("D473").Value = 45106
$ws.Range("E473").Value = 8
$ws.Range("F473").Value = 100112017
$ws.Range("G473").Value = "Apio"
$ws.Range("H473").Value = "Americana (o)"
$ws.Range("I473").Value = "Primera"
$ws.Range("J473").Value = 100
$ws.Range("K473").Value = 7000
$ws.Range("L473").Value = 7500
$ws.Range("M473").Value = 7250
$ws.Range("N473").Value = "$/docena de matas"
$ws.Range("O473").Value = "Región de Coquimbo"
$ws.Range("P473").Value = 1208
$ws.Range("Q473").Value = 6
$ws.Range("R473").Value = "Hortaliza"

# Populate the newly inserted row 474 (Segunda) with this week's data.
$ws.Range("A474").Value = 11
$ws.Range("B474").Value = "Vega Monumental Concepción"
$ws.Range("C474").Value = "Bíobío"
$ws.Range("D474").Value = 45106
$ws.Range("E474").Value = 8
$ws.Range("F474").Value = 100112017
$ws.Range("G474").Value = "Apio"
$ws.Range("H474").Value = "Americana (o)"
$ws.Range("I474").Value = "Segunda"
$ws.Range("J474").Value = 50
$ws.Range("K474").Value = 6500
$ws.Range("L474").Value = 6500
$ws.Range("M474").Value = 6500
$ws.Range("N474").Value = "$/docena de matas"
$ws.Range("O474").Value = "Región de Coquimbo"
$ws.Range("P474").Value = 1083
$ws.Range("Q474").Value = 6
$ws.Range("R474").Value = "Hortaliza"
